$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$c = $ws.Cells.Item(2,47)
$c.NumberFormat = "General"
